$wb = $excel.ActiveWorkbook

# ---- Metrics sheet: update the raw metric values (B2:B13) ----
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 337457.43999999994
$metrics.Range("B3").Value  = 272643.77
$metrics.Range("B4").Value  = 106545.01
$metrics.Range("B5").Value  = 13296
$metrics.Range("B6").Value  = 4256708.3199999994
$metrics.Range("B7").Value  = 3600171.2499999995
$metrics.Range("B8").Value  = 1235910.6900000002
$metrics.Range("B9").Value  = 164456
$metrics.Range("B10").Value = 32722032.120999824
$metrics.Range("B11").Value = 19630041.320000004
$metrics.Range("B12").Value = 11517619.58
$metrics.Range("B13").Value = 1262083

# Move the remembered selection on Metrics to match the saved view state
$metrics.Range("E38").Select()

# ---- today sheet: the B11:B22/E11:E22/F11:F22 formulas recalc automatically ----
# from the Metrics edits above (they reference Metrics!B2:B13). Just restore
# the saved selection for this sheet's view state.
$today = $wb.Worksheets.Item("today")
$today.Range("D5").Select()
